# StagingTemplates\Staging.OutcomeOrganization.xlsx
#
# The commit simply renames the three "key" columns on the staging sheet
# (the workbook was moved into the StagingTemplates directory and the
# column headers were relabeled to the new naming convention):
#
#   A2: OutcomeOrganization_ID   -> OrganizationBusinessKey
#   B2: OrganizationSourceKey    -> OutcomeBusinessKey
#   C2: OutcomeSourceKey         -> OutcomeOrganization_ID
#
# Row 1 (the "For internal use only..." banner) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "OrganizationBusinessKey"
$ws.Range("B2").Value = "OutcomeBusinessKey"
$ws.Range("C2").Value = "OutcomeOrganization_ID"
